$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: rename existing push notification string and add new email notification column
$ws.Range("D5").Value = "up_RidePushUserNotification 10"
$ws.Range("E5").Value = "up_RideEmailUserNotification 10"

# New row 17: "ride like" event
$ws.Range("A17").Value = "ride like"
$ws.Range("B17").Value = "owner"
$ws.Range("C17").Value = "ride detail"
$ws.Range("D17").Value = "up_NotificationFriendEmailDevice 1, 1"
$ws.Range("E17").Value = "up_NotificationFriendEmailDevice 1, 0"

$ws.Range("C17:E17").Style = $ws.Range("C16").Style
